# Change the title of this presentation (slide 1's title placeholder).
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$s.Shapes.Title.TextFrame.TextRange.Text = "Kanban vs. Scrum"
